$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Re-style row 20: it is no longer the last row, so it becomes a standalone
# 1-row group (top+bottom border, style indices 8/9) just like row 15. Formats only ---
$ws.Range("A15:E15").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ===== Scene: rows [21] =====
$ws.Range("A15:E15").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(21).RowHeight = 43.2
$ws.Range("C21").Value = ' We''ll get it done next time!'
$ws.Range("A21").Value = 'SCRIPT/D06P11A/um0802.ssb'
$ws.Range("D21").Value = ' В следующий раз у нас всё\nполучится!'
$ws.Range("E21").Value = ' Â òìåäôýþéê ñàè ô îàò âòæ\nðïìôœéóòÿ!'
$ws.Range("B21").Value = 681

# ===== Scene: rows [22, 23] =====
$ws.Range("A16:E16").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("A19:E19").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Range("C22").Value = ' All we gotta do now is cross\nover this here mountain. Then we''ll be at base\ncamp! Yup yup!'
$ws.Range("C23").Value = ' Let''s do it!'
$ws.Range("A22").Value = 'SCRIPT/D07P11A/um0802.ssb'
$ws.Range("D22").Value = ' Всё что нам нужно, это пересечь\nгору. И мы окажемся в лагере! Да-да!'
$ws.Range("D23").Value = ' Сделаем это!'
$ws.Range("E22").Value = ' Âòæ œóï îàí îôçîï, üóï ðåñåòåœû\nãïñô. É íú ïëàçåíòÿ â ìàãåñå! Äà-äà!'
$ws.Range("E23").Value = ' Òäåìàåí üóï!'
$ws.Range("B22").Value = 659
$ws.Range("B23").Value = 662

# ===== Scene: rows [24, 25] =====
$ws.Range("A16:E16").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)
$ws.Range("A19:E19").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(24).RowHeight = 43.2
$ws.Range("C24").Value = ' We just have to cross this here\nmountain! That''s all! Then we''ll be at the base\ncamp, yup yup!'
$ws.Range("C25").Value = ' Let''s do it!'
$ws.Range("A24").Value = 'SCRIPT/D07P11A/um0803.ssb'
$ws.Range("D24").Value = ' Нам всего лишь нужно пересечь\nэту гору! И всё! Мы окажемся в лагере!\nДа-да!'
$ws.Range("D25").Value = ' Сделаем это!'
$ws.Range("E24").Value = ' Îàí âòåãï ìéšû îôçîï ðåñåòåœû\nüóô ãïñô! É âòæ! Íú ïëàçåíòÿ â ìàãåñå!\nÄà-äà!'
$ws.Range("E25").Value = ' Òäåìàåí üóï!'
$ws.Range("B24").Value = 637
$ws.Range("B25").Value = 640

# ===== Scene: rows [26, 27] =====
$ws.Range("A16:E16").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)
$ws.Range("A16:E16").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(26).RowHeight = 43.2
$ws.Range("C26").Value = ' Oooh...[K] The great [CS:N]Dusknoir[CR] shook\nmy hand! By gosh and by golly…'
$ws.Range("C27").Value = ' It''s like...[K] Like a dream...[K] Wow…'
$ws.Range("A26").Value = 'SCRIPT/T01P01A/um1106.ssb'
$ws.Range("D26").Value = ' Оооо...[K] Великий [CS:N]Даскнуар[CR] пожал\nмне лапу! Вот это да, ей богу...'
$ws.Range("D27").Value = ' Я словно...[K] Во сне...[K] Вау...'
$ws.Range("E26").Value = ' Ïïïï...[K] Âåìéëéê [CS:N]Äàòëîôàñ[CR] ðïçàì\níîå ìàðô! Âïó üóï äà, åê áïãô...'
$ws.Range("E27").Value = ' Ÿ òìïâîï...[K] Âï òîå...[K] Âàô...'
$ws.Range("B26").Value = 612
$ws.Range("B27").Value = 616

# --- Row 27 must NOT have a column-A cell at all (middle-of-group-like row, no filename,
# matching rows 4/5/10/13's pattern where A is entirely absent). ---
$ws.Range("A27").ClearContents()
$ws.Range("A27").ClearFormats()

# --- Match the author's final on-screen state: scrolled down and C26 selected ---
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
[void]$ws.Range("C26").Select()

Write-Host "Done adding rows 21-27"